$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the status / timestamp cells for the "35efe67b..." row (row 2)
#    on all three sheets: Overview, zh-cn, de-de.
# ---------------------------------------------------------------------------

$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B2").Value = "Ready for handoff"
$ovw.Range("C2").Value = "Ready for handoff"
$ovw.Range("D2").Value = "2016-37-20 18:37:24"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-20 18:37:21"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-20 18:37:24"

# ---------------------------------------------------------------------------
# 2. Remove the "86b20423..." row (row 3) from all three sheets - that file
#    is no longer part of the handoff report.
# ---------------------------------------------------------------------------

$ovw.Rows.Item(3).Delete()
$zhcn.Rows.Item(3).Delete()
$dede.Rows.Item(3).Delete()
